$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.7302023333333333
$ws.Range("H2").Value = 2.190607
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1278363333333334
$ws.Range("N2").Value = 0.383509
$ws.Range("O2").Value = 0.002480915078704262
$ws.Range("P2").Value = 0.002480915078704262
$ws.Range("Q2").Value = 0.0933463888847778
$ws.Range("R2").Value = 0.8401174999630001
$ws.Range("S2").Value = 0.002480915078704262
$ws.Range("T2").Value = 0.002480915078704262

# Row 3
$ws.Range("G3").Value = 0.7302023333333333
$ws.Range("H3").Value = 2.190607
$ws.Range("O3").Value = 0.002269935507489869
$ws.Range("P3").Value = 0.002269935507489869
$ws.Range("Q3").Value = 0.08540811591833332
$ws.Range("R3").Value = 0.7686730432649999
$ws.Range("S3").Value = 0.002269935507489869
$ws.Range("T3").Value = 0.002269935507489869

# Row 4
$ws.Range("G4").Value = 0.7302023333333333
$ws.Range("H4").Value = 2.190607
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04120633333333334
$ws.Range("N4").Value = 0.123619
$ws.Range("O4").Value = 0.0007996898146180199
$ws.Range("P4").Value = 0.0007996898146180199
$ws.Range("Q4").Value = 0.03008896074811111
$ws.Range("R4").Value = 0.270800646733
$ws.Range("S4").Value = 0.0007996898146180199
$ws.Range("T4").Value = 0.0007996898146180199

# Row 5
$ws.Range("G5").Value = 0.7302023333333333
$ws.Range("H5").Value = 2.190607
$ws.Range("M5").Value = 51.241888
$ws.Range("N5").Value = 153.725664
$ws.Range("O5").Value = 0.9944494595991877
$ws.Range("P5").Value = 0.9944494595991878
$ws.Range("Q5").Value = 37.41694618200533
$ws.Range("R5").Value = 336.752515638048
$ws.Range("S5").Value = 0.9944494595991877
$ws.Range("T5").Value = 0.9944494595991878
